# Word COM-interop script implementing the diff described in the task.
$d = $word.ActiveDocument
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARN: replace failed for: $find"
    }
}

function Find-ParagraphIndex($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

# --- A. Simple text replacements (paragraph count unaffected) ---

Replace-Text "a basic festival area with stores which consists" "a basic festival area which consists"

Replace-Text "Guest, Stores and an Information Centre" "Guest, Store and an Information Centre"

Replace-Text "Guests: Simulating the virtue of humans at a fest, Guests possess the qualities of a human as feeling hungry, thirsty, getting food and drink.  " "Guest: Simulating the virtue of a human at a festival, Guest possesses the qualities of a human as feeling hungry, thirsty, getting food and drink and moving.  "

Replace-Text "the final model servers as a main file" "the final model serves as a main file"

Replace-Text "Using functions such as Reflexes and Skills we demonstrate the behavior of agents." "Using functions such as Reflex and Skill, we demonstrate the behavior of agents."

Replace-Text "Each species will have “Reflexes” to act" "Each species will have “Reflex” to act"

Replace-Text "Reflex beCrazy: This is when guests wander around the fest." "Reflex beCrazy: when a guest wanders around the fest."

Replace-Text "Reflex thirstyOrHungry: This to let the guest know that they are hungry or thirsty" "Reflex thirstyOrHungry: let the guest know that they are hungry or thirsty"

Replace-Text "Reflex moveToTarget: Guests go to a specified target" "Reflex moveToTarget: Guest goes to a specified target"

Replace-Text "Reflex askInfoCentre: Guests go to Information centre to fetch information" "Reflex askInfoCentre: Guest goes to Information centre to fetch information"

Replace-Text "Reflex drinkOrEat: Guests fulfill their hunger or thirst." "Reflex drinkOrEat: Guest fulfills his hunger or thirst."

Replace-Text "As they would want to explore different shops, they can ask the nearby guests. This is implemented" "They can either directly use their memory to go to a shop or, as they want to discover new shops, they can go to the Information centre or ask the nearby guests on the way. This is implemented"

Replace-Text "Reflex askGuest: Introduces the logic" "Reflex askGuest: introduces the logic"

Replace-Text "gave us a significant knowledge" "gave us significant knowledge"

Replace-Text "to simulate the multi-agents systems" "to simulate multi-agents systems"

Replace-Text "GAMA platform and it’s syntax" "GAMA platform and its syntax"

# --- B. "Agents" section paragraph restructure: Shops-paragraph becomes the new
#     Store-paragraph (losing its indent override) and the old Information Centre
#     paragraph's slot now holds the new Information Centre text (formatting kept). ---

$storeSlot = Find-ParagraphIndex "Information Centre: Simulating the virtue of the Information Desk at a fest which has most of the information about the fest in general.*"
if ($storeSlot -eq -1) { Write-Output "WARN: could not locate Information Centre paragraph (store slot)" }
$xmlStore = "<w:p $ns><w:pPr><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val='0'/></w:rPr><w:t xml:space='preserve'>Store: Simulating the virtue of a store, which can sell food and drink. In this simulation, there are two attributes of a shop, either it sells food or drink whenever a guest asks for it.</w:t></w:r></w:p>"
$d.Paragraphs($storeSlot).Range.InsertXML($xmlStore)

$infoSlot = Find-ParagraphIndex "Shops: Simulating the virtue of the shops*"
if ($infoSlot -eq -1) { Write-Output "WARN: could not locate Shops paragraph (info slot)" }
$p = $d.Paragraphs($infoSlot)
$r = $p.Range
$rTrim = $d.Range($r.Start, $r.End - 1)
$rTrim.Text = "Information Centre: Simulating the virtue of the Information Desk at a festival which has most of the information about it in general. "

# --- C. New empty paragraphs before the "Model 2" and "Model 3" headings,
#     plus formatting changes around the "Model 3" heading. ---

$model2idx = Find-ParagraphIndex "Model 2:*"
if ($model2idx -eq -1) { Write-Output "WARN: could not locate Model 2 heading" }
$d.Paragraphs($model2idx).Range.InsertParagraphBefore()
$xmlCenterEmpty = "<w:p $ns><w:pPr><w:ind w:left='0' w:firstLine='0'/><w:jc w:val='center'/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val='0'/></w:rPr></w:r></w:p>"
$d.Paragraphs($model2idx).Range.InsertXML($xmlCenterEmpty)

$model3idx = Find-ParagraphIndex "Model 3:*"
if ($model3idx -eq -1) { Write-Output "WARN: could not locate Model 3 heading" }
$d.Paragraphs($model3idx).Range.InsertParagraphBefore()
$d.Paragraphs($model3idx).Range.InsertXML($xmlCenterEmpty)

$model3idx = $model3idx + 1
$xmlModel3Heading = "<w:p $ns><w:pPr><w:ind w:left='0' w:firstLine='0'/><w:jc w:val='left'/><w:rPr><w:color w:val='4a86e8'/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val='4a86e8'/><w:rtl w:val='0'/></w:rPr><w:t xml:space='preserve'>Model 3: Challenge 2: Festival_security.gaml</w:t></w:r></w:p>"
$d.Paragraphs($model3idx).Range.InsertXML($xmlModel3Heading)

$d.Paragraphs($model3idx).Range.InsertParagraphAfter()
$xmlModel3TrailingEmpty = "<w:p $ns><w:pPr><w:ind w:left='0' w:firstLine='0'/><w:jc w:val='left'/><w:rPr><w:color w:val='4a86e8'/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val='0'/></w:rPr></w:r></w:p>"
$d.Paragraphs($model3idx + 1).Range.InsertXML($xmlModel3TrailingEmpty)

Write-Output "Done"
